$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New row 10 values: A10 = 4 (hours), B10 = new topic text
$ws.Range("A10").Value = 4
$ws.Range("B10").Value = "Implemented exp, pow, div, sub, Neuron, Layers"

# B9 (existing "Automatic Backpropogation..." row) gets wrap-text only alignment
$ws.Range("B9").WrapText = $true

# Update the active selection to match the post-edit state (B15)
$ws.Range("B15").Select()
